$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header: column E label changes from "micRef" to "cRef"
$ws.Range("E1").Value = "cRef"

# Row 2 (amoxicillin): stock 32 -> 4, cRef 8 -> 2, Solvent DMSO -> H2O
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = "H2O"

# Row 4 (colistin): cRef 1 -> 2
$ws.Range("E4").Value = 2

# Row 5 (fosfomycin): cRef 1 -> 4
$ws.Range("E5").Value = 4

# Update selection to match final state
$ws.Range("F3").Select()
